$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove John McGuire from the mailing list (row 27 of the "Trans Page" sheet)
$ws.Rows(27).Delete()

# Roll the commission period forward from October 2017 to November 2017:
# update the subject line text, then every attachment file name / path.
[void]$ws.Cells.Replace("October 2017", "November 2017")
[void]$ws.Cells.Replace("201710", "201711")
